$d = $word.ActiveDocument

$replacements = @(
    @{old="2024-05-20 Monday"; new="2024-05-21 Tuesday"},
    @{old="430÷8="; new="386÷4="},
    @{old="152÷3="; new="904÷7="},
    @{old="314÷6="; new="860÷4="},
    @{old="715÷2="; new="894÷5="},
    @{old="442÷6="; new="503÷3="},
    @{old="378÷7="; new="395÷8="},
    @{old="931÷9="; new="229÷8="},
    @{old="981÷6="; new="621÷9="},
    @{old="684÷8="; new="235÷7="},
    @{old="464÷7="; new="258÷5="},
    @{old="573÷7="; new="399÷3="},
    @{old="554÷4="; new="652÷5="},
    @{old="497÷9="; new="722÷5="},
    @{old="930÷2="; new="658÷3="},
    @{old="643÷3="; new="747÷2="},
    @{old="745÷6="; new="212÷6="},
    @{old="495÷3="; new="136÷8="},
    @{old="404÷7="; new="908÷5="},
    @{old="900÷6="; new="546÷8="},
    @{old="267÷6="; new="437÷5="},
    @{old="581÷9="; new="340÷6="},
    @{old="489÷9="; new="578÷7="},
    @{old="691÷6="; new="443÷3="},
    @{old="547÷6="; new="502÷9="},
    @{old="702÷3="; new="657÷8="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
